$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - TestCase_004 (NASA APOD service)
$ws.Range("B5").Value = "https://api.nasa.gov/planetary/apod?api_key=NNKOjkoul8n1CH18TWA9gwngW1s1SmjESPjNoUFo"
$ws.Range("C5").Value = "media_type"
$ws.Range("D5").Value = "image"
$ws.Range("A5").Value = "TestCase_004"
$ws.Range("E5").Value = "'200"

# Row 6 - TestCase_005 (NASA Mars Rover photos service)
$ws.Range("B6").Value = "https://api.nasa.gov/mars-photos/api/v1/rovers/curiosity/photos?sol=1000&page=2&api_key=DEMO_KEY"
$ws.Range("C6").Value = "photos[1].rover.name"
$ws.Range("D6").Value = "Curiosity"
$ws.Range("A6").Value = "TestCase_005"
$ws.Range("E6").Value = "'200"

$ws.Range("A7").Select()
